# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.927.01"
$ws.Range("E2").Value = "  +7.34%  "

$ws.Range("D3").Value = "3.487.49"
$ws.Range("E3").Value = "  +5.45%  "

$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "'415.25"
$ws.Range("E5").Value = "  +3.79%  "

$ws.Range("D6").Value = "'128.69"
$ws.Range("E6").Value = "  +17.08%  "

$ws.Range("D7").Value = "3.480.18"
$ws.Range("E7").Value = "  +5.39%  "

$ws.Range("D8").Value = "'0.593"
$ws.Range("E8").Value = "  +1.46%  "

$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("E10").Value = "  +8.88%  "

$ws.Range("E11").Value = "  +29.66%  "

$ws.Range("D12").Value = "'42.34"
$ws.Range("E12").Value = "  +6.17%  "

$ws.Range("E13").Value = "  +0.59%  "

$ws.Range("D14").Value = "4.038.62"
$ws.Range("E14").Value = "  +5.95%  "

$ws.Range("D15").Value = "'8.73"
$ws.Range("E15").Value = "  +4.10%  "

$ws.Range("D16").Value = "'20.06"
$ws.Range("E16").Value = "  +4.48%  "

$ws.Range("D17").Value = "3.506.14"
$ws.Range("E17").Value = "  +6.51%  "

$ws.Range("D18").Value = "62.932.25"
$ws.Range("E18").Value = "  +8.03%  "

$ws.Range("D19").Value = "'1.04"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("D20").Value = "'10.82"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("D21").Value = "'0.0000136"
$ws.Range("E21").Value = "  +24.24%  "

$ws.Range("D22").Value = "'3.36"
$ws.Range("E22").Value = "  +0.76%  "

$ws.Range("D23").Value = "'81.98"
$ws.Range("E23").Value = "  +9.62%  "

$ws.Range("D24").Value = "'314.65"
$ws.Range("E24").Value = "  +4.02%  "

$ws.Range("D25").Value = "'13.09"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").Value = "'30.84"
$ws.Range("E27").Value = "  +8.65%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'8.06"
$ws.Range("E28").Value = "  +2.40%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.74"
$ws.Range("E29").Value = "  +4.69%  "

$ws.Range("D30").Value = "'0.179"
$ws.Range("E30").Value = "  +5.22%  "

$ws.Range("D31").Value = "'4.36"
$ws.Range("E31").Value = "  -1.26%  "

$ws.Range("D32").Value = "'0.118"
$ws.Range("E32").Value = "  +4.45%  "

$ws.Range("D33").Value = "'2.64"
$ws.Range("E33").Value = "  +23.43%  "

$ws.Range("D34").Value = "'11.70"
$ws.Range("E34").Value = "  +2.75%  "

$ws.Range("D35").Value = "'42.93"
$ws.Range("E35").Value = "  +4.34%  "

$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("E37").Value = "  -6.72%  "

$ws.Range("D38").Value = "'52.40"
$ws.Range("E38").Value = "  +1.08%  "

$ws.Range("D39").Value = "'3.56"
$ws.Range("E39").Value = "  +1.93%  "

$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("E41").Value = "  -7.54%  "

$ws.Range("D42").Value = "'2.02"
$ws.Range("E42").Value = "  +7.18%  "

$ws.Range("E43").Value = "  +2.30%  "

$ws.Range("D44").Value = "'136.19"
$ws.Range("E44").Value = "  -1.25%  "

$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D45").Value = "'17.15"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.287"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("E47").Value = "  -0.17%  "

$ws.Range("E48").Value = "  -1.85%  "

$ws.Range("D49").Value = "'21.97"
$ws.Range("E49").Value = "  -2.53%  "

$ws.Range("D50").Value = "2.229.72"
$ws.Range("E50").Value = "  +2.72%  "

$ws.Range("D51").Value = "3.835.52"
$ws.Range("E51").Value = "  +5.80%  "
